$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("H2").Value = 3.1
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63
$ws.Range("Q2").Value = 2.35
$ws.Range("R2").Value = 1.57
$ws.Range("S2").Value = 1.53
$ws.Range("T2").Value = 2.38
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 1.73
$ws.Range("AA2").Value = 21
$ws.Range("AC2").Value = 7
$ws.Range("AF2").Value = 67
$ws.Range("AK2").Value = 34
$ws.Range("AN2").Value = 13
$ws.Range("AR2").Value = 2.38
$ws.Range("AS2").Value = 9
$ws.Range("AZ2").Value = 301
$ws.Range("BB2").Value = 251

# Row 3 updates
$ws.Range("G3").Value = 9.5
$ws.Range("H3").Value = 5.25
$ws.Range("J3").Value = 9.5
$ws.Range("L3").Value = 1.8
